# "add participant to groupe": update the recipient list in column A.
# Row 1 keeps the "email" header; row 2 gets a new participant's address;
# row 3's old mailto hyperlink/contact is replaced by a new plain address.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "droitformation.web@gmail.com"
$ws.Range("A3").Value = "hello@yahoo.fr"

# The old A3 value carried a mailto hyperlink + the built-in "hyperlink"
# cell style - drop both so the cell goes back to plain, unformatted text.
$ws.Hyperlinks.Delete()
$wb.Styles("Lien hypertexte").Delete()
$ws.Range("A3").ClearFormats()
